# edit.ps1 -- applies "updating kickoff doc for partners" changes
# to Guest Authentication Through Our Partners.docx

$d = $word.ActiveDocument

$ENDASH = [char]0x2013
$RSQUOTE = [char]0x2019

function Replace-InRange($range, [string]$old, [string]$new) {
    $ok = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Replace failed: old=[$old] new=[$new]"
    }
}

# ---------------------------------------------------------------------------
# 1. "Please forward the sample OIDC or SAML metadata ..." -- turn the word
#    "SAML" into a hyperlink pointing at the SAML metadata URL.
# ---------------------------------------------------------------------------
$pFwd = $d.Paragraphs(10).Range
$pFwdStart = $pFwd.Start
$pFwdText = $pFwd.Text
$samlIdx = $pFwdText.IndexOf("SAML metadata")
$samlStart = $pFwdStart + $samlIdx
$samlEnd = $samlStart + 4
$samlRange = $d.Range($samlStart, $samlEnd)
$d.Hyperlinks.Add($samlRange, "https://alaska-poc.cic-demo-platform.auth0app.com/samlp/metadata/DPHf8btcMeuYWFVSWnZIh0Q41gryQmQG") | Out-Null

Write-Output "Step 1 done"

# ---------------------------------------------------------------------------
# 2. "Create the Partner application in our Identity Provider" gains
#    " (Auth0)" and lower-cases "provider".
# ---------------------------------------------------------------------------
Replace-InRange $d.Paragraphs(13).Range `
    "Create the Partner application in our Identity Provider" `
    "Create the Partner application in our Identity provider (Auth0)"

Write-Output "Step 2 done"

# ---------------------------------------------------------------------------
# 3. Scope-of-work bullet restructure:
#    - remove "Email client credentials ... TEST environment" bullet (ilvl2)
#    - replace "Configure Partner redirect URLs ..." bullet (ilvl2) with
#      "Configure Partner's SAML ACS url (or OIDC redirect url)"
#    - add a new ilvl1 bullet "Send Alaska's SSO login url ..."
#    - update the "Additional work may be needed ..." bullet's wording
# ---------------------------------------------------------------------------

# remove the "Email client credentials" bullet (paragraph 14) completely
$d.Paragraphs(14).Range.Delete()

# paragraph 14 is now "Configure Partner redirect URLs..." -- replace it
Replace-InRange $d.Paragraphs(14).Range `
    "Configure Partner redirect URLs for SAML or OIDC integrations in our Identity platform" `
    "Configure Partner$($RSQUOTE)s SAML ACS url (or OIDC redirect url)"

# insert a new bullet after it, one level up (ilvl1), with the SSO login text
$d.Paragraphs(14).Range.InsertParagraphAfter()
$newBullet = $d.Paragraphs(15)
$newBullet.Range.ListFormat.ListLevelNumber = 2
$newBullet.Range.Text = "Send Alaska$($RSQUOTE)s SSO login url with client id to Partner so that they can begin integration testing against our TEST environment"

# update the "Additional work may be needed ..." bullet (now paragraph 16)
Replace-InRange $d.Paragraphs(16).Range `
    "business logic which is tightly coupled with Ping to return the Guest$($RSQUOTE)s Identity data points" `
    "business logic which is tightly coupled with Ping to return additional Guest$($RSQUOTE)s Identity data points"

Write-Output "Step 3 done"

# ---------------------------------------------------------------------------
# 4. Timeline section.
#    - split "9/27/24 ... provides client id to Partner" into two bullets:
#        "9/27/24 - Auth0 provisions TEST environment for AS"
#        "10/7/24 - Alaska Airlines creates applications in TEST environment
#         for each Partner and provides SSO login uri with client id to Partner"
#    - "10/4/24" becomes "10/11/24"
#    - add a new bold/red "TBD QA environment" paragraph
#    - "10/18/24 ... provides client id to Partner" gains "SSO login uri with"
# ---------------------------------------------------------------------------

Replace-InRange $d.Paragraphs(49).Range `
    "9/27/24 $($ENDASH) Alaska Airlines creates applications in TEST environment for each Partner and provides client id to Partner" `
    "10/7/24 $($ENDASH) Alaska Airlines creates applications in TEST environment for each Partner and provides SSO login uri with client id to Partner"

$d.Paragraphs(49).Range.InsertParagraphBefore()
$newTimelinePara = $d.Paragraphs(49)
$newTimelinePara.Range.Text = "9/27/24 $($ENDASH) Auth0 provisions TEST environment for AS"

# "10/4/24" -> "10/11/24" (paragraph shifted down by 1, now 51)
Replace-InRange $d.Paragraphs(51).Range "10/4/24" "10/11/24"

# new bold/red "TBD QA environment" paragraph, inserted after the 10/11/24 bullet
$d.Paragraphs(51).Range.InsertParagraphAfter()
$tbdPara = $d.Paragraphs(52)
$tbdPara.Range.Text = "TBD QA environment"
$tbdPara.Range.Font.Bold = 1
$tbdPara.Range.Font.BoldBi = 1
$tbdPara.Range.Font.Color = 255  ## wdColorRed (BGR 0x0000FF == RGB FF0000)

# "10/18/24 ... provides client id to Partner" (now paragraph 53)
Replace-InRange $d.Paragraphs(53).Range `
    "creates applications in PROD environment for each Partner and provides client id to Partner" `
    "creates applications in PROD environment for each Partner and provides SSO login uri with client id to Partner"

Write-Output "Step 4 done"
